# Re-orders data rows 3-26 of the sheet according to the mapping below.
# Mapping key = destination row, value = source row (both referring to the
# original, pre-edit layout of the worksheet). This reproduces a re-sort of
# the observation rows while leaving the header rows (1-2) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 26
$lastCol = "AY"

# destination row -> source row (1-based sheet row numbers)
$rowMap = @{
    3  = 15
    4  = 3
    5  = 5
    6  = 12
    7  = 16
    8  = 10
    9  = 13
    10 = 26
    11 = 22
    12 = 14
    13 = 17
    14 = 9
    15 = 6
    16 = 21
    17 = 18
    18 = 19
    19 = 25
    20 = 20
    21 = 24
    22 = 11
    23 = 8
    24 = 23
    25 = 7
    26 = 4
}

# Columns that contain plain date-looking text (e.g. "2023-08-26"). When such
# text is written through Value2 into a "General" formatted cell, Excel
# auto-detects it as a real date and stores it as a date serial number with a
# date number format instead of keeping the original text. To avoid that, we
# temporarily force those columns to a text format before assigning values,
# and restore their original style afterwards so the resulting cells look
# exactly like the untouched ones (no left-over explicit style).
$dateTextCols = @("Y", "AA")
$origDateStyles = @{}
foreach ($col in $dateTextCols) {
    $origDateStyles[$col] = $ws.Range($col + "1").Style
    $ws.Range($col + $firstRow + ":" + $col + $lastRow).NumberFormat = "@"
}

# Read the whole block once so we can freely rearrange the rows without any
# cell being overwritten before it has been read.
$srcRange = $ws.Range("A" + $firstRow + ":" + $lastCol + $lastRow)
$srcValues = $srcRange.Value2

$numRows = $lastRow - $firstRow + 1
$numCols = $srcValues.GetLength(1)

$newValues = New-Object 'object[,]' $numRows, $numCols

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    $destIdx = $destRow - $firstRow
    $srcIdx = $srcRow - $firstRow + 1
    for ($c = 1; $c -le $numCols; $c++) {
        $newValues[$destIdx, ($c - 1)] = $srcValues[$srcIdx, $c]
    }
}

$destRange = $ws.Range("A" + $firstRow + ":" + $lastCol + $lastRow)
$destRange.Value2 = $newValues

foreach ($col in $dateTextCols) {
    $ws.Range($col + $firstRow + ":" + $col + $lastRow).Style = $origDateStyles[$col]
}
